# Correção nos dados e início da análise PNAD 2009
#
# The sheet had two "category header" rows that only carried a label
# (no B/C/D values): row 5 "situação do domicílio" and row 8 "grandes
# regiões e unidades da federação". Those rows are removed entirely so
# that the "urbana"/"rural" and "norte"/"rondônia"/... rows move up and
# sit directly under "brasil", which also corrects the row->data
# alignment for every region/state below them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "situação do domicílio" header row.
$ws.Rows.Item(5).Delete()

# After the row above was removed, the "grandes regiões e unidades da
# federação" header row (originally row 8) is now row 7.
$ws.Rows.Item(7).Delete()

Write-Host "Final UsedRange:" $ws.UsedRange.Address()
